$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the label text for row 3 (A3)
$ws.Range("A3").Value = "2020-06-29_diff"

# Update the numeric values for row 3 (B3:D3)
$ws.Range("B3").Value = 0.2822615868947058
$ws.Range("C3").Value = 2.239176668751999
$ws.Range("D3").Value = -10.28507807938547
